# Update countries & provincias Spain
#
# This script applies the 30-Apr-2020 00:52 data refresh to the "Pais"
# worksheet:
#   - Refreshes the "Datos actualizados..." timestamp in A1.
#   - Updates case/death/recovered counters for several countries whose
#     numbers changed between the 00:22 and 00:52 snapshots.
#   - Colombia overtakes Panama in the ranking (rows 49/50 swap country,
#     keep their sorted position).
#   - Guinea-Bisau's case count jumped sharply, moving it from rank 167
#     up to rank 138; every country that used to occupy ranks 138-162
#     shifts down by one rank (row 134 is filled with Guinea-Bisau's new
#     totals, row 163 now holds the country that used to be directly
#     above Guinea-Bisau, i.e. Islas Caimanes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 30 de Abril de 2020 a las 00:52"

# Row 4
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1057152
$ws.Cells.Item(4, 3).Value = 21387
$ws.Cells.Item(4, 4).Value = 145389
$ws.Cells.Item(4, 5).Value = 850414
$ws.Cells.Item(4, 6).Value = 18665
$ws.Cells.Item(4, 7).Value = 2083
$ws.Cells.Item(4, 8).Value = 61349

# Row 10
$ws.Cells.Item(10, 1).Value = "Turquia"
$ws.Cells.Item(10, 2).Value = 117589
$ws.Cells.Item(10, 3).Value = 2936
$ws.Cells.Item(10, 4).Value = 44040
$ws.Cells.Item(10, 5).Value = 70468
$ws.Cells.Item(10, 6).Value = 1574
$ws.Cells.Item(10, 7).Value = 89
$ws.Cells.Item(10, 8).Value = 3081

# Row 45
$ws.Cells.Item(45, 1).Value = "Chequia"
$ws.Cells.Item(45, 2).Value = 7579
$ws.Cells.Item(45, 3).Value = 75
$ws.Cells.Item(45, 4).Value = 3108
$ws.Cells.Item(45, 5).Value = 4244
$ws.Cells.Item(45, 6).Value = 71
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 227

# Row 49
$ws.Cells.Item(49, 1).Value = "Colombia"
$ws.Cells.Item(49, 2).Value = 6207
$ws.Cells.Item(49, 3).Value = 258
$ws.Cells.Item(49, 4).Value = 1411
$ws.Cells.Item(49, 5).Value = 4518
$ws.Cells.Item(49, 6).Value = 118
$ws.Cells.Item(49, 7).Value = 9
$ws.Cells.Item(49, 8).Value = 278

# Row 50
$ws.Cells.Item(50, 1).Value = "Panama"
$ws.Cells.Item(50, 2).Value = 6200
$ws.Cells.Item(50, 3).Value = 179
$ws.Cells.Item(50, 4).Value = 484
$ws.Cells.Item(50, 5).Value = 5540
$ws.Cells.Item(50, 6).Value = 89
$ws.Cells.Item(50, 7).Value = 9
$ws.Cells.Item(50, 8).Value = 176

# Row 118
$ws.Cells.Item(118, 1).Value = "Kenia"
$ws.Cells.Item(118, 2).Value = 384
$ws.Cells.Item(118, 3).Value = 10
$ws.Cells.Item(118, 4).Value = 129
$ws.Cells.Item(118, 5).Value = 240
$ws.Cells.Item(118, 6).Value = 2
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 15

# Row 134
$ws.Cells.Item(134, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(134, 2).Value = 205
$ws.Cells.Item(134, 3).Value = 132
$ws.Cells.Item(134, 4).Value = 19
$ws.Cells.Item(134, 5).Value = 185
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 1

# Row 135
$ws.Cells.Item(135, 1).Value = "Islas Feroe"
$ws.Cells.Item(135, 2).Value = 187
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 181
$ws.Cells.Item(135, 5).Value = 6
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 0

# Row 136
$ws.Cells.Item(136, 1).Value = "Martinica"
$ws.Cells.Item(136, 2).Value = 175
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 83
$ws.Cells.Item(136, 5).Value = 78
$ws.Cells.Item(136, 6).Value = 5
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 14

# Row 137
$ws.Cells.Item(137, 1).Value = "Birmania"
$ws.Cells.Item(137, 2).Value = 150
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 27
$ws.Cells.Item(137, 5).Value = 117
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 1
$ws.Cells.Item(137, 8).Value = 6

# Row 138
$ws.Cells.Item(138, 1).Value = "Guadalupe"
$ws.Cells.Item(138, 2).Value = 149
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 82
$ws.Cells.Item(138, 5).Value = 55
$ws.Cells.Item(138, 6).Value = 11
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 12

# Row 139
$ws.Cells.Item(139, 1).Value = "Liberia"
$ws.Cells.Item(139, 2).Value = 141
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 45
$ws.Cells.Item(139, 5).Value = 80
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 16

# Row 140
$ws.Cells.Item(140, 1).Value = "Gibraltar"
$ws.Cells.Item(140, 2).Value = 141
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(140, 4).Value = 131
$ws.Cells.Item(140, 5).Value = 10
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 0

# Row 141
$ws.Cells.Item(141, 1).Value = "Brunei"
$ws.Cells.Item(141, 2).Value = 138
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 124
$ws.Cells.Item(141, 5).Value = 13
$ws.Cells.Item(141, 6).Value = 2
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 1

# Row 142
$ws.Cells.Item(142, 1).Value = "Etiopia"
$ws.Cells.Item(142, 2).Value = 130
$ws.Cells.Item(142, 3).Value = 4
$ws.Cells.Item(142, 4).Value = 58
$ws.Cells.Item(142, 5).Value = 69
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 3

# Row 143
$ws.Cells.Item(143, 1).Value = "Madagascar"
$ws.Cells.Item(143, 2).Value = 128
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 90
$ws.Cells.Item(143, 5).Value = 38
$ws.Cells.Item(143, 6).Value = 1
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 0

# Row 144
$ws.Cells.Item(144, 1).Value = "Guayana Francesa"
$ws.Cells.Item(144, 2).Value = 125
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 93
$ws.Cells.Item(144, 5).Value = 31
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 1

# Row 145
$ws.Cells.Item(145, 1).Value = "Camboya"
$ws.Cells.Item(145, 2).Value = 122
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 119
$ws.Cells.Item(145, 5).Value = 3
$ws.Cells.Item(145, 6).Value = 1
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 0

# Row 146
$ws.Cells.Item(146, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(146, 2).Value = 116
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 68
$ws.Cells.Item(146, 5).Value = 40
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 8

# Row 147
$ws.Cells.Item(147, 1).Value = "Cabo Verde"
$ws.Cells.Item(147, 2).Value = 114
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 2
$ws.Cells.Item(147, 5).Value = 111
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 1

# Row 148
$ws.Cells.Item(148, 1).Value = "Bermudas"
$ws.Cells.Item(148, 2).Value = 111
$ws.Cells.Item(148, 3).Value = 1
$ws.Cells.Item(148, 4).Value = 48
$ws.Cells.Item(148, 5).Value = 57
$ws.Cells.Item(148, 6).Value = 10
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 6

# Row 149
$ws.Cells.Item(149, 1).Value = "Togo"
$ws.Cells.Item(149, 2).Value = 109
$ws.Cells.Item(149, 3).Value = 10
$ws.Cells.Item(149, 4).Value = 64
$ws.Cells.Item(149, 5).Value = 38
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 7

# Row 150
$ws.Cells.Item(150, 1).Value = "Sierra Leona"
$ws.Cells.Item(150, 2).Value = 104
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 12
$ws.Cells.Item(150, 5).Value = 88
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 4

# Row 151
$ws.Cells.Item(151, 1).Value = "Aruba"
$ws.Cells.Item(151, 2).Value = 100
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 73
$ws.Cells.Item(151, 5).Value = 25
$ws.Cells.Item(151, 6).Value = 4
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 2

# Row 152
$ws.Cells.Item(152, 1).Value = "Zambia"
$ws.Cells.Item(152, 2).Value = 97
$ws.Cells.Item(152, 3).Value = 2
$ws.Cells.Item(152, 4).Value = 54
$ws.Cells.Item(152, 5).Value = 40
$ws.Cells.Item(152, 6).Value = 1
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 3

# Row 153
$ws.Cells.Item(153, 1).Value = "Monaco"
$ws.Cells.Item(153, 2).Value = 95
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 58
$ws.Cells.Item(153, 5).Value = 33
$ws.Cells.Item(153, 6).Value = 1
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 4

# Row 154
$ws.Cells.Item(154, 1).Value = "Suazilandia"
$ws.Cells.Item(154, 2).Value = 91
$ws.Cells.Item(154, 3).Value = 20
$ws.Cells.Item(154, 4).Value = 10
$ws.Cells.Item(154, 5).Value = 80
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 1

# Row 155
$ws.Cells.Item(155, 1).Value = "Liechtenstein"
$ws.Cells.Item(155, 2).Value = 82
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 55
$ws.Cells.Item(155, 5).Value = 26
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 1

# Row 156
$ws.Cells.Item(156, 1).Value = "Uganda"
$ws.Cells.Item(156, 2).Value = 81
$ws.Cells.Item(156, 3).Value = 2
$ws.Cells.Item(156, 4).Value = 52
$ws.Cells.Item(156, 5).Value = 29
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 0

# Row 157
$ws.Cells.Item(157, 1).Value = "Bahamas"
$ws.Cells.Item(157, 2).Value = 80
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 23
$ws.Cells.Item(157, 5).Value = 46
$ws.Cells.Item(157, 6).Value = 1
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 11

# Row 158
$ws.Cells.Item(158, 1).Value = "Barbados"
$ws.Cells.Item(158, 2).Value = 80
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 39
$ws.Cells.Item(158, 5).Value = 34
$ws.Cells.Item(158, 6).Value = 4
$ws.Cells.Item(158, 7).Value = 1
$ws.Cells.Item(158, 8).Value = 7

# Row 159
$ws.Cells.Item(159, 1).Value = "Mozambique"
$ws.Cells.Item(159, 2).Value = 76
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = 12
$ws.Cells.Item(159, 5).Value = 64
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

# Row 160
$ws.Cells.Item(160, 1).Value = "Haiti"
$ws.Cells.Item(160, 2).Value = 76
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 8
$ws.Cells.Item(160, 5).Value = 62
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 6

# Row 161
$ws.Cells.Item(161, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(161, 2).Value = 75
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 33
$ws.Cells.Item(161, 5).Value = 29
$ws.Cells.Item(161, 6).Value = 7
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 13

# Row 162
$ws.Cells.Item(162, 1).Value = "Guyana"
$ws.Cells.Item(162, 2).Value = 74
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 15
$ws.Cells.Item(162, 5).Value = 51
$ws.Cells.Item(162, 6).Value = 5
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 8

# Row 163
$ws.Cells.Item(163, 1).Value = "Islas Caimanes"
$ws.Cells.Item(163, 2).Value = 73
$ws.Cells.Item(163, 3).Value = 3
$ws.Cells.Item(163, 4).Value = 10
$ws.Cells.Item(163, 5).Value = 62
$ws.Cells.Item(163, 6).Value = 3
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 1

# Row 204
$ws.Cells.Item(204, 1).Value = "Surinam"
$ws.Cells.Item(204, 2).Value = 10
$ws.Cells.Item(204, 3).Value = 0
$ws.Cells.Item(204, 4).Value = 8
$ws.Cells.Item(204, 5).Value = 1
$ws.Cells.Item(204, 6).Value = 0
$ws.Cells.Item(204, 7).Value = 0
$ws.Cells.Item(204, 8).Value = 1
